$d = $word.ActiveDocument

function Set-BoldPrefix($paraIndex, $prefixText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $found = $r.Find.Execute($prefixText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Bold = $true
    }
}

# Задача 1.
Set-BoldPrefix 2 "Задача 1."
# Решение. (paragraph 4)
Set-BoldPrefix 4 "Решение."
# Ответ: 35.
Set-BoldPrefix 9 "Ответ:"
# Задача 2.
Set-BoldPrefix 10 "Задача 2."
# Решение.  (paragraph 12, bold prefix includes trailing space)
Set-BoldPrefix 12 "Решение. "
# Ответ: 48.  (bold prefix includes trailing space)
Set-BoldPrefix 17 "Ответ: "
# Задача 3.
Set-BoldPrefix 18 "Задача 3."
# Решение.  (paragraph 19) - whole paragraph becomes bold (incl. paragraph mark)
$p19 = $d.Paragraphs(19)
$p19.Range.Bold = $true
# Ответ: 102.
Set-BoldPrefix 21 "Ответ: "
# Задача 4.
Set-BoldPrefix 22 "Задача 4."
# Решение. (paragraph 24)
Set-BoldPrefix 24 "Решение."
# Ответ: 6.
Set-BoldPrefix 38 "Ответ: "
# Задача 5.
Set-BoldPrefix 39 "Задача 5."
# Решение. (paragraph 41)
Set-BoldPrefix 41 "Решение."
# Ответ: 3.
Set-BoldPrefix 65 "Ответ:"

# Домашнее задание heading: center + 14pt (sz/szCs = 28 half-points)
$p67 = $d.Paragraphs(67)
$p67.Alignment = 1
$p67.Range.Font.Size = 14
$p67.Range.Font.SizeBi = 14

Write-Output "edits applied"
